# New Microsite scripts support to Beta server
# Appends new sprint-run rows to the AMSIN sheet (rows 20-23) and the AMS
# sheet (rows 10-11), matching the upstream UI_EDUCATION_REGISTRATION_HISTORY_DATA
# data refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AMSIN": append rows 20..23
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Column A holds plain text dates (e.g. "2022-09-15"), not real Excel date
# serials. Writing the string straight into .Value lets the host
# auto-recognise the ISO-ish text as a date and coerce it to a date serial,
# which is not what the source workbook stores. Routing the literal text
# through a formula ("="..."") and then collapsing the formula down to a
# static value with Copy/PasteSpecial(xlPasteValues) keeps it as real text.
$wsAmsin.Cells.Item(20, 1).Formula = '="2022-09-15"'
$wsAmsin.Cells.Item(21, 1).Formula = '="2022-09-16"'
$wsAmsin.Cells.Item(22, 1).Formula = '="2022-09-19"'
$wsAmsin.Cells.Item(23, 1).Formula = '="2022-09-20"'
$wsAmsin.Range("A20:A23").Copy()
$wsAmsin.Range("A20:A23").PasteSpecial(-4163)  # xlPasteValues

# Column B is a real numeric date/time serial - copy the existing row's
# number format down first so new rows keep the "YYYY-MM-DD HH:MM:SS" style.
$wsAmsin.Range("B19").Copy()
$wsAmsin.Range("B20:B23").PasteSpecial(-4122)  # xlPasteFormats

$wsAmsin.Cells.Item(20, 2).Value = 44819.62414569444
$wsAmsin.Cells.Item(21, 2).Value = 44820.64496002315
$wsAmsin.Cells.Item(22, 2).Value = 44823.60046166667
$wsAmsin.Cells.Item(23, 2).Value = 44824.38674135417

# Column C (sprint/script name), D/E/F (case counts) and G (time taken) are
# plain text/numbers - no coercion risk.
$wsAmsin.Cells.Item(20, 3).Value = "eduecs166"
$wsAmsin.Cells.Item(21, 3).Value = "fstcedu167"
$wsAmsin.Cells.Item(22, 3).Value = "scndedu167"
$wsAmsin.Cells.Item(23, 3).Value = "finaleduc167"

$wsAmsin.Cells.Item(20, 4).Value = 60
$wsAmsin.Cells.Item(21, 4).Value = 60
$wsAmsin.Cells.Item(22, 4).Value = 60
$wsAmsin.Cells.Item(23, 4).Value = 60

$wsAmsin.Cells.Item(20, 5).Value = 60
$wsAmsin.Cells.Item(21, 5).Value = 60
$wsAmsin.Cells.Item(22, 5).Value = 60
$wsAmsin.Cells.Item(23, 5).Value = 60

$wsAmsin.Cells.Item(20, 6).Value = 0
$wsAmsin.Cells.Item(21, 6).Value = 0
$wsAmsin.Cells.Item(22, 6).Value = 0
$wsAmsin.Cells.Item(23, 6).Value = 0

$wsAmsin.Cells.Item(20, 7).Value = 1.33
$wsAmsin.Cells.Item(21, 7).Value = 1.23
$wsAmsin.Cells.Item(22, 7).Value = 2.95
$wsAmsin.Cells.Item(23, 7).Value = 1.17

# ---------------------------------------------------------------------------
# Sheet "AMS": row 10 gets restyled/re-synced, and a new row 11 is appended
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Bring row 10's formatting in line with the rest of the table (A/C/D/E/F/G
# pick up the same look as row 9) while leaving column B's existing
# date/time number format (already applied) untouched.
$wsAms.Range("A9").Copy()
$wsAms.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$wsAms.Range("C9:G9").Copy()
$wsAms.Range("C10:G10").PasteSpecial(-4122)  # xlPasteFormats

# Re-sync row 10's run time to the refreshed source value.
$wsAms.Cells.Item(10, 2).Value = 44812.52689822917

# New row 11: same text-coercion guard as above for the plain-text run date.
$wsAms.Cells.Item(11, 1).Formula = '="2022-09-20"'
$wsAms.Range("A11").Copy()
$wsAms.Range("A11").PasteSpecial(-4163)  # xlPasteValues

$wsAms.Range("B10").Copy()
$wsAms.Range("B11").PasteSpecial(-4122)  # xlPasteFormats
$wsAms.Cells.Item(11, 2).Value = 44824.72934250939

$wsAms.Cells.Item(11, 3).Value = "betaedu167"
$wsAms.Cells.Item(11, 4).Value = 60
$wsAms.Cells.Item(11, 5).Value = 60
$wsAms.Cells.Item(11, 6).Value = 0
$wsAms.Cells.Item(11, 7).Value = 0.83
